# Update the underlying "nb" counts in the three tables. All dependent
# formulas (D, E, F, G, H columns and the various SUM() totals) recalc
# automatically because they are formulas in the sheet already.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Table 1 (rows 4-10) ---
$ws.Range("B5").Value = 133633
$ws.Range("B6").Value = 39314
$ws.Range("B7").Value = 5166
$ws.Range("B8").Value = 1662
$ws.Range("C8").Value = 0
$ws.Range("B9").Value = 3078
$ws.Range("C9").Value = 1

# H10 now sums one more row (H4:H8 instead of H4:H7)
$ws.Range("H10").Formula = "=SUM(H4:H8)"

# --- Table 2 (rows 12-18) ---
$ws.Range("B13").Value = 133633
$ws.Range("B14").Value = 39312
$ws.Range("C14").Value = 2
$ws.Range("B15").Value = 5164
$ws.Range("C15").Value = 3
$ws.Range("C16").Value = 7
$ws.Range("C17").Value = 1

# --- Table 3 (rows 21-27) ---
$ws.Range("B22").Value = 133633
$ws.Range("B23").Value = 39314
$ws.Range("C23").Value = 0
$ws.Range("B24").Value = 5167
$ws.Range("C24").Value = 1507
$ws.Range("B25").Value = 127
$ws.Range("C25").Value = 1535
$ws.Range("B26").Value = 3078
$ws.Range("C26").Value = 28

$excel.Calculate()
